# Apply "Add data for 2022-08-01" update to the carjacking-by-neighborhood
# workbook: rename the sheet/header from "July 23" to "July 24" and update
# the per-neighborhood monthly counts in column B (and a handful of other
# historical columns) to reflect the newly-added day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab and update the column header label in B1.
$ws.Name = "Through 2022-07-24"
$ws.Range("B1").Value = "July 2022 (through July 24)"

# Helper-less direct cell updates (row number -> neighborhood per column A).
$ws.Range("B2").Value  = 12   # Austin
$ws.Range("P2").Value  = 7    # Austin

$ws.Range("I3").Value  = 5    # Englewood
$ws.Range("AY3").Value = 2    # Englewood

$ws.Range("P5").Value  = 13   # Garfield Park

$ws.Range("AY7").Value = 1    # Roseland (new value)

$ws.Range("B8").Value  = 7    # North Lawndale
$ws.Range("P8").Value  = 18   # North Lawndale

$ws.Range("W13").Value = 1    # River North (new value)

$ws.Range("AR14").Value = 1   # West Pullman (new value)

$ws.Range("B16").Value = 3    # Douglas

$ws.Range("AK24").Value = 1   # South Deering (new value)

$ws.Range("B26").Value = 5    # Little Village

$ws.Range("B29").Value = 7    # Humboldt Park

$ws.Range("P35").Value = 2    # United Center

$ws.Range("I38").Value = 6    # West Town

$ws.Range("B41").Value = 3    # Logan Square
$ws.Range("AK41").Value = 2   # Logan Square

$ws.Range("B47").Value = 2    # Little Italy, UIC

$ws.Range("W48").Value = 2    # Irving Park

$ws.Range("AK52").Value = 3   # Chatham

$ws.Range("I57").Value = 2    # Woodlawn

$ws.Range("I61").Value = 1    # Avalon Park (new value)

$ws.Range("P66").Value = 2    # Chinatown
